# ITST-16721 IPP xls2csv: wrong handling of custom format - adjust format
# matching, so that the custom format does correctly get replaced - adapt
# test.
#
# Adds a new row (row 6) to the worksheet that exercises a custom number
# format: B6 holds the numeric value 268.02 formatted with the custom
# format "[$-10409]#,##0.00;\-#,##0.00" (right aligned / top aligned),
# and C6 holds an explanatory text describing the scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B6: numeric value using the new custom number format.
$b6 = $ws.Range("B6")
$b6.Value = 268.02
$b6.NumberFormat = "[$-10409]#,##0.00;\-#,##0.00"
$b6.HorizontalAlignment = -4152   # xlRight
$b6.VerticalAlignment = -4160     # xlTop
$b6.ReadingOrder = 1              # xlRightToLeft

# C6: explanation text (new shared string).
$ws.Range("C6").Value = "Format should be interpreted correctly."
